$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last populated row in column B (data rows)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row

# Header for the new column C
$ws.Range("C1").Value = "Superficie_m2"

# Fill formulas for each data row: Superficie_m2 = Superficie km2 * 1000 * 1000
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Formula = "=B$r*1000*1000"
}

# Autosize the new column to fit its contents
$ws.Columns.Item(3).EntireColumn.AutoFit() | Out-Null

# Match the selection left behind in the saved file
$ws.Range("E7").Select() | Out-Null
